$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add formula D2 = C2*B2
$ws.Range("D2").Formula = "=C2*B2"

# Update selection to match final state (F5)
$ws.Range("F5").Select()

# Update window position
$excel.ActiveWindow.Left = 180
$excel.ActiveWindow.Top = 1530
